# Apply the edit described by the commit:
# "implemented first working function and samples generating valid GTFS feed"
#
# Concretely, on sheet "Tabelle1":
#   - Cell E5 text changes from "BfFestAlbtal" to "BahnhofsfestAlbtal"
#   - The active selection moves from E1 (whole-column sqref) to the single cell E5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the shared-string cell content.
$ws.Range("E5").Value = "BahnhofsfestAlbtal"

# Move/collapse the selection onto E5, matching the new cursor position.
$ws.Activate()
$ws.Range("E5").Select() | Out-Null
